# "Generate Report for Handoff"
#
# Status flips from "In Translation" to "Ready for handoff" (every sheet
# that showed that status), and the two timestamps that record when the
# handoff xliff was generated move forward a few dozen seconds. The
# "Status" / "zh-cn" / "de-de" columns also get a bit wider so the new,
# longer "Ready for handoff" label isn't truncated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Latest HO Xliff Generate Date / Latest Handoff Datetime -------------
$wsOverview.Range("G2").Value = "2016-08-31 19:07:25"
$wsDeDe.Range("H2").Value     = "2016-08-31 19:07:25"
$wsZhCn.Range("H2").Value     = "2016-08-31 19:07:21"

# --- Widen the Status / zh-cn / de-de columns to fit the new text --------
# (stored column width goes from ~13.41 chars to ~17.22 chars; ColumnWidth
# is pixel-quantized by Excel, so we request the value that lands closest
# to the target on that grid)
$wsOverview.Range("E1").ColumnWidth = 16.333333333333332
$wsOverview.Range("F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth     = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth     = 16.333333333333332
